# feat: add 2022-Q3 data
#
# 1. Insert a new "2022-Q3" worksheet right after "总计", containing the
#    fund-holding detail table for that quarter.
# 2. Update the "总计" (summary) sheet so it gains a new leading data row
#    for 2022-Q3 (all the other quarter rows shift down by one).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: add the new "2022-Q3" sheet right after "总计"
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

# Header row (bold, centered, top-aligned, thin border - matches the other
# quarter sheets in this workbook)
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"
$q3.Range("B1:H1").Font.Bold = $true
$q3.Range("B1:H1").HorizontalAlignment = -4108
$q3.Range("B1:H1").VerticalAlignment = -4160
$q3.Range("B1:H1").Borders.LineStyle = 1

# Data row 2: 014320 / 德邦半导体产业混合C
$q3.Range("A2").Value = 0
$q3.Range("B2").NumberFormat = "@"
$q3.Range("B2").Value = "014320"
$q3.Range("C2").Value = "德邦半导体产业混合C"
$q3.Range("D2").NumberFormat = "@"
$q3.Range("D2").Value = "1.38"
$q3.Range("E2").NumberFormat = "@"
$q3.Range("E2").Value = "91.65"
$q3.Range("F2").NumberFormat = "@"
$q3.Range("F2").Value = "6.04"
$q3.Range("G2").NumberFormat = "@"
$q3.Range("G2").Value = "0.0834"
$q3.Range("H2").Value = 3

# Data row 3: 014319 / 德邦半导体产业混合A
$q3.Range("A3").Value = 1
$q3.Range("B3").NumberFormat = "@"
$q3.Range("B3").Value = "014319"
$q3.Range("C3").Value = "德邦半导体产业混合A"
$q3.Range("D3").NumberFormat = "@"
$q3.Range("D3").Value = "0.37"
$q3.Range("E3").NumberFormat = "@"
$q3.Range("E3").Value = "91.65"
$q3.Range("F3").NumberFormat = "@"
$q3.Range("F3").Value = "6.04"
$q3.Range("G3").NumberFormat = "@"
$q3.Range("G3").Value = "0.0223"
$q3.Range("H3").Value = 3

$q3.Range("A2:A3").Font.Bold = $true
$q3.Range("A2:A3").HorizontalAlignment = -4108
$q3.Range("A2:A3").VerticalAlignment = -4160
$q3.Range("A2:A3").Borders.LineStyle = 1

# Match the page-margin convention used by every other sheet in this workbook
$q3.PageSetup.LeftMargin = 54
$q3.PageSetup.RightMargin = 54
$q3.PageSetup.TopMargin = 72
$q3.PageSetup.BottomMargin = 72
$q3.PageSetup.HeaderMargin = 36
$q3.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------
# Step 2: update the "总计" sheet - insert a new row 2 for 2022-Q3 and
# push the existing quarters (2022-Q2 ... 2021-Q1) down by one row
# ---------------------------------------------------------------------
$total.Rows("2:2").Insert()
$total.Range("A2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.11

$total.Range("A2").Font.Bold = $true
$total.Range("A2").HorizontalAlignment = -4108
$total.Range("A2").VerticalAlignment = -4160
$total.Range("A2").Borders.LineStyle = 1

# ---------------------------------------------------------------------
# Step 3: restore "总计" as the active tab (matches the unchanged bookViews
# in the source diff - only the sheet order/content changed)
# ---------------------------------------------------------------------
$total.Activate()
